# Restructure rows 7-8, 12-14, and 17-28: rows that were a lone section
# header ("Education Completed", "Household", "Age Migrated", "Migration
# Cohort", "Acculturation") swap places with data rows that followed them,
# per the target layout. We rewrite every touched row explicitly (values
# + clearing the B:L tail on rows that become header-only) so the result
# matches cell-for-cell regardless of the starting state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Less than Primary
$ws.Cells.Item(7, 1).Value = 'Less than Primary'
$ws.Cells.Item(7, 2).Value = '0.4'
$ws.Cells.Item(7, 3).Value = '0.17'
$ws.Cells.Item(7, 4).Value = '0.28'
$ws.Cells.Item(7, 5).Value = '0.12'
$ws.Cells.Item(7, 6).Value = '0.29'
$ws.Cells.Item(7, 7).Value = '0.1'
$ws.Cells.Item(7, 8).Value = '0.11'
$ws.Cells.Item(7, 9).Value = '0.07'
$ws.Cells.Item(7, 10).Value = '0.03'
$ws.Cells.Item(7, 11).Value = '0.01'
$ws.Cells.Item(7, 12).Value = '0.03'

# Row 8: Education Completed
$ws.Cells.Item(8, 1).Value = 'Education Completed'
$ws.Range("B8:L8").ClearContents()

# Row 12: Household Size
$ws.Cells.Item(12, 1).Value = 'Household Size'
$ws.Cells.Item(12, 2).Value = '3.45'
$ws.Cells.Item(12, 3).Value = '2.44'
$ws.Cells.Item(12, 4).Value = '2.93'
$ws.Cells.Item(12, 5).Value = '2.57'
$ws.Cells.Item(12, 6).Value = '3.24'
$ws.Cells.Item(12, 7).Value = '2.82'
$ws.Cells.Item(12, 8).Value = '2.86'
$ws.Cells.Item(12, 9).Value = '2.37'
$ws.Cells.Item(12, 10).Value = '2.06'
$ws.Cells.Item(12, 11).Value = '1.95'
$ws.Cells.Item(12, 12).Value = '2.27'

# Row 13: Lives Alone
$ws.Cells.Item(13, 1).Value = 'Lives Alone'
$ws.Cells.Item(13, 2).Value = '0.15'
$ws.Cells.Item(13, 3).Value = '0.3'
$ws.Cells.Item(13, 4).Value = '0.23'
$ws.Cells.Item(13, 5).Value = '0.27'
$ws.Cells.Item(13, 6).Value = '0.18'
$ws.Cells.Item(13, 7).Value = '0.2'
$ws.Cells.Item(13, 8).Value = '0.19'
$ws.Cells.Item(13, 9).Value = '0.28'
$ws.Cells.Item(13, 10).Value = '0.39'
$ws.Cells.Item(13, 11).Value = '0.29'
$ws.Cells.Item(13, 12).Value = '0.31'

# Row 14: Household
$ws.Cells.Item(14, 1).Value = 'Household'
$ws.Range("B14:L14").ClearContents()

# Row 17: 25 - 49
$ws.Cells.Item(17, 1).Value = '25 - 49'
$ws.Cells.Item(17, 2).Value = '0.77'
$ws.Cells.Item(17, 3).Value = '0.61'
$ws.Cells.Item(17, 4).Value = '0.69'
$ws.Cells.Item(17, 5).Value = '0.63'
$ws.Cells.Item(17, 6).Value = '0.79'
$ws.Cells.Item(17, 7).Value = '0.7'
$ws.Cells.Item(17, 8).Value = '0.7'
$ws.Cells.Item(17, 9).Value = '-'
$ws.Cells.Item(17, 10).Value = '-'
$ws.Cells.Item(17, 11).Value = '-'
$ws.Cells.Item(17, 12).Value = '-'

# Row 18: 50 and Above
$ws.Cells.Item(18, 1).Value = '50 and Above'
$ws.Cells.Item(18, 2).Value = '0.23'
$ws.Cells.Item(18, 3).Value = '0.39'
$ws.Cells.Item(18, 4).Value = '0.31'
$ws.Cells.Item(18, 5).Value = '0.37'
$ws.Cells.Item(18, 6).Value = '0.21'
$ws.Cells.Item(18, 7).Value = '0.3'
$ws.Cells.Item(18, 8).Value = '0.3'
$ws.Cells.Item(18, 9).Value = '1'
$ws.Cells.Item(18, 10).Value = '1'
$ws.Cells.Item(18, 11).Value = '1'
$ws.Cells.Item(18, 12).Value = '1'

# Row 19: Before 1965
$ws.Cells.Item(19, 1).Value = 'Before 1965'
$ws.Cells.Item(19, 2).Value = '0.03'
$ws.Cells.Item(19, 3).Value = '0.05'
$ws.Cells.Item(19, 4).Value = '0.03'
$ws.Cells.Item(19, 5).Value = '0.08'
$ws.Cells.Item(19, 6).Value = '0.02'
$ws.Cells.Item(19, 7).Value = '0.03'
$ws.Cells.Item(19, 8).Value = '0.04'
$ws.Cells.Item(19, 9).Value = '-'
$ws.Cells.Item(19, 10).Value = '-'
$ws.Cells.Item(19, 11).Value = '-'
$ws.Cells.Item(19, 12).Value = '-'

# Row 20: Age Migrated
$ws.Cells.Item(20, 1).Value = 'Age Migrated'
$ws.Range("B20:L20").ClearContents()

# Row 21: 1965 - 1979
$ws.Cells.Item(21, 1).Value = '1965 - 1979'
$ws.Cells.Item(21, 2).Value = '0.29'
$ws.Cells.Item(21, 3).Value = '0.22'
$ws.Cells.Item(21, 4).Value = '0.2'
$ws.Cells.Item(21, 5).Value = '0.34'
$ws.Cells.Item(21, 6).Value = '0.24'
$ws.Cells.Item(21, 7).Value = '0.23'
$ws.Cells.Item(21, 8).Value = '0.25'
$ws.Cells.Item(21, 9).Value = '-'
$ws.Cells.Item(21, 10).Value = '-'
$ws.Cells.Item(21, 11).Value = '-'
$ws.Cells.Item(21, 12).Value = '-'

# Row 22: 1980 - 1999
$ws.Cells.Item(22, 1).Value = '1980 - 1999'
$ws.Cells.Item(22, 2).Value = '0.52'
$ws.Cells.Item(22, 3).Value = '0.43'
$ws.Cells.Item(22, 4).Value = '0.55'
$ws.Cells.Item(22, 5).Value = '0.38'
$ws.Cells.Item(22, 6).Value = '0.61'
$ws.Cells.Item(22, 7).Value = '0.47'
$ws.Cells.Item(22, 8).Value = '0.5'
$ws.Cells.Item(22, 9).Value = '-'
$ws.Cells.Item(22, 10).Value = '-'
$ws.Cells.Item(22, 11).Value = '-'
$ws.Cells.Item(22, 12).Value = '-'

# Row 23: After 1999
$ws.Cells.Item(23, 1).Value = 'After 1999'
$ws.Cells.Item(23, 2).Value = '0.2'
$ws.Cells.Item(23, 3).Value = '0.35'
$ws.Cells.Item(23, 4).Value = '0.27'
$ws.Cells.Item(23, 5).Value = '0.32'
$ws.Cells.Item(23, 6).Value = '0.19'
$ws.Cells.Item(23, 7).Value = '0.31'
$ws.Cells.Item(23, 8).Value = '0.25'
$ws.Cells.Item(23, 9).Value = '-'
$ws.Cells.Item(23, 10).Value = '-'
$ws.Cells.Item(23, 11).Value = '-'
$ws.Cells.Item(23, 12).Value = '-'

# Row 24: Migration Cohort
$ws.Cells.Item(24, 1).Value = 'Migration Cohort'
$ws.Range("B24:L24").ClearContents()

# Row 25: Citizen
$ws.Cells.Item(25, 1).Value = 'Citizen'
$ws.Cells.Item(25, 2).Value = '0.45'
$ws.Cells.Item(25, 3).Value = '-'
$ws.Cells.Item(25, 4).Value = '0.64'
$ws.Cells.Item(25, 5).Value = '0.75'
$ws.Cells.Item(25, 6).Value = '0.61'
$ws.Cells.Item(25, 7).Value = '0.67'
$ws.Cells.Item(25, 8).Value = '0.73'
$ws.Cells.Item(25, 9).Value = '-'
$ws.Cells.Item(25, 10).Value = '-'
$ws.Cells.Item(25, 11).Value = '-'
$ws.Cells.Item(25, 12).Value = '-'

# Row 26: English Speakers
$ws.Cells.Item(26, 1).Value = 'English Speakers'
$ws.Cells.Item(26, 2).Value = '0.64'
$ws.Cells.Item(26, 3).Value = '0.84'
$ws.Cells.Item(26, 4).Value = '0.62'
$ws.Cells.Item(26, 5).Value = '0.66'
$ws.Cells.Item(26, 6).Value = '0.77'
$ws.Cells.Item(26, 7).Value = '0.84'
$ws.Cells.Item(26, 8).Value = '0.9'
$ws.Cells.Item(26, 9).Value = '0.99'
$ws.Cells.Item(26, 10).Value = '1'
$ws.Cells.Item(26, 11).Value = '1'
$ws.Cells.Item(26, 12).Value = '1'

# Row 27: N
$ws.Cells.Item(27, 1).Value = 'N'
$ws.Cells.Item(27, 2).Value = '44152'
$ws.Cells.Item(27, 3).Value = '9159'
$ws.Cells.Item(27, 4).Value = '6541'
$ws.Cells.Item(27, 5).Value = '13647'
$ws.Cells.Item(27, 6).Value = '13655'
$ws.Cells.Item(27, 7).Value = '17935'
$ws.Cells.Item(27, 8).Value = '227703'
$ws.Cells.Item(27, 9).Value = '120724'
$ws.Cells.Item(27, 10).Value = '313063'
$ws.Cells.Item(27, 11).Value = '3165675'
$ws.Cells.Item(27, 12).Value = '94162'

# Row 28: Acculturation
$ws.Cells.Item(28, 1).Value = 'Acculturation'
$ws.Range("B28:L28").ClearContents()
